$d = $word.ActiveDocument

$d.Content.Find.Execute("287×6=", $true, $false, $false, $false, $false, $true, 1, $false, "364×8=", 2) | Out-Null
$d.Content.Find.Execute("523×8=", $true, $false, $false, $false, $false, $true, 1, $false, "269×5=", 2) | Out-Null
$d.Content.Find.Execute("318×8=", $true, $false, $false, $false, $false, $true, 1, $false, "351×2=", 2) | Out-Null
$d.Content.Find.Execute("947×9=", $true, $false, $false, $false, $false, $true, 1, $false, "459×3=", 2) | Out-Null
$d.Content.Find.Execute("200×3=", $true, $false, $false, $false, $false, $true, 1, $false, "421×8=", 2) | Out-Null
$d.Content.Find.Execute("693×6=", $true, $false, $false, $false, $false, $true, 1, $false, "326×2=", 2) | Out-Null
$d.Content.Find.Execute("113×9=", $true, $false, $false, $false, $false, $true, 1, $false, "711×4=", 2) | Out-Null
$d.Content.Find.Execute("670×8=", $true, $false, $false, $false, $false, $true, 1, $false, "795×3=", 2) | Out-Null
$d.Content.Find.Execute("607×6=", $true, $false, $false, $false, $false, $true, 1, $false, "837×4=", 2) | Out-Null
$d.Content.Find.Execute("686×5=", $true, $false, $false, $false, $false, $true, 1, $false, "718×7=", 2) | Out-Null
$d.Content.Find.Execute("852×9=", $true, $false, $false, $false, $false, $true, 1, $false, "627×4=", 2) | Out-Null
$d.Content.Find.Execute("559×4=", $true, $false, $false, $false, $false, $true, 1, $false, "885×3=", 2) | Out-Null
$d.Content.Find.Execute("730×9=", $true, $false, $false, $false, $false, $true, 1, $false, "676×8=", 2) | Out-Null
$d.Content.Find.Execute("554×5=", $true, $false, $false, $false, $false, $true, 1, $false, "889×8=", 2) | Out-Null
$d.Content.Find.Execute("310×9=", $true, $false, $false, $false, $false, $true, 1, $false, "877×4=", 2) | Out-Null
$d.Content.Find.Execute("133×8=", $true, $false, $false, $false, $false, $true, 1, $false, "595×5=", 2) | Out-Null
$d.Content.Find.Execute("450×7=", $true, $false, $false, $false, $false, $true, 1, $false, "511×4=", 2) | Out-Null
$d.Content.Find.Execute("129×5=", $true, $false, $false, $false, $false, $true, 1, $false, "800×7=", 2) | Out-Null
$d.Content.Find.Execute("509×4=", $true, $false, $false, $false, $false, $true, 1, $false, "829×3=", 2) | Out-Null
$d.Content.Find.Execute("736×4=", $true, $false, $false, $false, $false, $true, 1, $false, "639×7=", 2) | Out-Null
$d.Content.Find.Execute("412×3=", $true, $false, $false, $false, $false, $true, 1, $false, "958×2=", 2) | Out-Null
$d.Content.Find.Execute("425×6=", $true, $false, $false, $false, $false, $true, 1, $false, "853×9=", 2) | Out-Null
$d.Content.Find.Execute("916×7=", $true, $false, $false, $false, $false, $true, 1, $false, "658×2=", 2) | Out-Null
$d.Content.Find.Execute("638×5=", $true, $false, $false, $false, $false, $true, 1, $false, "329×4=", 2) | Out-Null
$d.Content.Find.Execute("102×3=", $true, $false, $false, $false, $false, $true, 1, $false, "717×8=", 2) | Out-Null

Write-Output "Replacements complete"
